$wb = $excel.ActiveWorkbook

# --- Sheet1: selection moved from K8 to C11 (tab focus also moves away from Sheet1) ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("C11").Select()

# --- Duplicate Sheet3 (keeps its original formatting/defaults) and insert the copy
#     right before it; the copy becomes the new "Sheet4" ---
$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet3.Copy($sheet3)
$newSheet = $wb.Worksheets.Item("Sheet3 (2)")
$newSheet.Name = "Sheet4"

# --- Populate "Sheet4" with the table data (rows 6-16, columns A/C/D) ---
for ($r = 6; $r -le 16; $r++) {
    $newSheet.Cells.Item($r, 1).Value = "sheet4"
    $newSheet.Cells.Item($r, 3).Value = "B COL "
    $newSheet.Cells.Item($r, 4).Value = "C COL"
}

# Column D width customization (closest achievable to the target stored width of 19.5546875,
# the engine quantizes ColumnWidth to 1/6-character steps)
$newSheet.Columns.Item(4).ColumnWidth = 18.666666666666668

# Select entire first row and make this new sheet the active tab
$newSheet.Range("A1:XFD1").Select()
$newSheet.Activate()
